$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Biological discard" row (row 40) and the "Shells NA" row (row 49).
# Delete from the bottom up so row indices of earlier rows remain valid.
$ws.Rows.Item(49).Delete()
$ws.Rows.Item(40).Delete()

# After the deletions the data block runs from row 35 to row 50.
# Set the W(kg) column (G) to 0 for every remaining data row in that block.
$ws.Range("G35:G50").Value = 0

# Anadara transversa (now row 36) also had its Numb (H) value changed to -1.
$ws.Range("H36").Value = -1
